$d = $word.ActiveDocument

# 1. Remove the trailing " AAAA" after "Noctambus" in the Mandant line.
$d.Content.Find.Execute("Noctambus AAAA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Noctambus", 2)

# 2. Rename the "Pardeliste" (List Paragraph) style back to "Paragraphedeliste".
$style = $d.Styles("Pardeliste")
$style.NameLocal = "Paragraphedeliste"
